$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "User: שלמה_בונצל`nEmail: shlezi0@gmail.com`nAction: Updated task`nTask type: ראיון מועמד לחונכות`nChanges: 1 field(s)`n  • Status: 'הושלמה' → 'בביצוע'`nRoles: System Administrator`nRecord ID: 189"
$ws.Range("A3").Value = "User: שלמה_בונצל`nEmail: shlezi0@gmail.com`nAction: Deleted pending tutor`nVolunteer: בהככ כדגד`nReason: Promoted to Tutor`nTask ID: 189`nStatus: Successfully promoted`nRoles: System Administrator"
$ws.Range("A4").Value = "User: שלמה_בונצל`nEmail: shlezi0@gmail.com`nAction: Updated task`nTask type: ראיון מועמד לחונכות`nChanges: 1 field(s)`n  • Status: 'בביצוע' → 'לא הושלמה'`nRoles: System Administrator`nRecord ID: 189"

$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()
